$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141, pushing the existing rows 141-161 down to 142-162.
$ws.Rows.Item(141).Insert()

# Seed the new row with the same constant columns as its neighbour (row 142,
# which is the former row 141), then overwrite the columns that actually
# carry new data for this weekly entry.
$ws.Range("A141:R141").Value2 = $ws.Range("A142:R142").Value2

$ws.Range("D141").Value2 = 44504
$ws.Range("I141").Value2 = "Primera"
$ws.Range("J141").Value2 = 350
$ws.Range("K141").Value2 = 7500
$ws.Range("L141").Value2 = 8000
$ws.Range("M141").Value2 = 7786
$ws.Range("P141").Value2 = 1298
